# Update automatico via Actualizar 05-18-2020 07-14-13
# Adds a new row (17/5/2020) to the "Condicion_Pacientes" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data table (ListObject) currently spans A1:F48 - grow it by one row.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Populate the three known columns for the new row (49): Fecha, Pruebas
# Realizadas, Pruebas Positivas.
$ws.Range("A49").Value = "17/5/2020"
$ws.Range("B49").Value = 1214
$ws.Range("C49").Value = 149

# Move the active selection to C50, matching the author's next data-entry
# position, without disturbing the current scroll position.
$win = $ws.Application.ActiveWindow
$savedScrollRow = $win.ScrollRow
$savedScrollColumn = $win.ScrollColumn
[void]$ws.Range("C50").Select()
$win.ScrollRow = $savedScrollRow
$win.ScrollColumn = $savedScrollColumn
